$wb = $excel.ActiveWorkbook
$wsTest = $wb.Worksheets.Item("Test")
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTest)
$ws2.Name = "Configuration"

$ws2.Range("A2").Value = "Browser"
$ws2.Range("A1").Value = "Property"
$ws2.Range("B1").Value = "Value"
$ws2.Range("B2").Value = "Chrome"

$ws2.Columns.Item(1).ColumnWidth = 20.67
$ws2.Columns.Item(2).ColumnWidth = 29.33

$ws2.PageSetup.Orientation = 1

$ws2.Range("B2").Select() | Out-Null
